$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02406506557682775
$ws.Range("D2").Value = 0.1668829158416827
$ws.Range("E2").Value = 0.1607398286185102
$ws.Range("F2").Value = 1.636872026552723
$ws.Range("G2").Value = 1.002614734140138
$ws.Range("H2").Value = 1.0185464478825
$ws.Range("I2").Value = 0.8687822641251515
$ws.Range("J2").Value = 0.1928632368882006
$ws.Range("K2").Value = 2.054569687855121
$ws.Range("C3").Value = 0.02244693183169488
$ws.Range("D3").Value = 0.1612439996455777
$ws.Range("E3").Value = 0.1566321673715727
$ws.Range("F3").Value = 1.644518109295433
$ws.Range("G3").Value = 1.010280372564281
$ws.Range("H3").Value = 1.0295335993668
$ws.Range("I3").Value = 0.8739126022889891
$ws.Range("J3").Value = 0.1889610989209913
$ws.Range("K3").Value = 1.84682674065283
$ws.Range("C4").Value = 0.02144792046223643
$ws.Range("D4").Value = 0.1578141065464393
$ws.Range("E4").Value = 0.1541827072036064
$ws.Range("F4").Value = 1.65055622308968
$ws.Range("G4").Value = 1.016034338535931
$ws.Range("H4").Value = 1.037016044151045
$ws.Range("I4").Value = 0.8778641947310817
$ws.Range("J4").Value = 0.1866880443344598
$ws.Range("K4").Value = 1.71928890023662
$ws.Range("C5").Value = 0.0210394588206313
$ws.Range("D5").Value = 0.1564246357730497
$ws.Range("E5").Value = 0.1532028087965358
$ws.Range("F5").Value = 1.653353566480618
$ws.Range("G5").Value = 1.018641255609452
$ws.Range("H5").Value = 1.040249929316133
$ws.Range("I5").Value = 0.8796753011676657
$ws.Range("J5").Value = 0.1857925700123886
$ws.Range("K5").Value = 1.667322366121539
$ws.Range("C6").Value = 0.02097155273430218
$ws.Range("D6").Value = 0.156194415340849
$ws.Range("E6").Value = 0.153041201459633
$ws.Range("F6").Value = 1.653838372387952
$ws.Range("G6").Value = 1.019089928861604
$ws.Range("H6").Value = 1.040798058396774
$ws.Range("I6").Value = 0.8799881397919478
$ws.Range("J6").Value = 0.1856457358427193
$ws.Range("K6").Value = 1.658693791336759
$ws.Range("C7").Value = 0.02144241727214791
$ws.Range("D7").Value = 0.1577953341850105
$ws.Range("E7").Value = 0.154169417941258
$ws.Range("F7").Value = 1.650592586892728
$ws.Range("G7").Value = 1.016068436585698
$ws.Range("H7").Value = 1.037058910062512
$ws.Range("I7").Value = 0.8778878079185972
$ws.Range("J7").Value = 0.1866758429997688
$ws.Range("K7").Value = 1.718588033891535
$ws.Range("C8").Value = 0.02350827939955025
$ws.Range("D8").Value = 0.164931936510115
$ws.Range("E8").Value = 0.159308427633654
$ws.Range("F8").Value = 1.639228944477509
$ws.Range("G8").Value = 1.005039765874429
$ws.Range("H8").Value = 1.022181749207434
$ws.Range("I8").Value = 0.8703843886714608
$ws.Range("J8").Value = 0.19149222483518
$ws.Range("K8").Value = 1.982937386471519
$ws.Range("C9").Value = 0.02751539173902273
$ws.Range("D9").Value = 0.1791811814171922
$ws.Range("E9").Value = 0.1699630535359162
$ws.Range("F9").Value = 1.627652786528472
$ws.Range("G9").Value = 0.9917770574836879
$ws.Range("H9").Value = 0.9988688530028043
$ws.Range("I9").Value = 0.8620649339170043
$ws.Range("J9").Value = 0.2019166465732383
$ws.Range("K9").Value = 2.501409148267783
$ws.Range("C10").Value = 0.03043204459057591
$ws.Range("D10").Value = 0.189802579072861
$ws.Range("E10").Value = 0.1781446012021632
$ws.Range("F10").Value = 1.625745525089997
$ws.Range("G10").Value = 0.9872096164910857
$ws.Range("H10").Value = 0.9853406628631518
$ws.Range("I10").Value = 0.8599004281669593
$ws.Range("J10").Value = 0.2101800547291504
$ws.Range("K10").Value = 2.882352775953962
$ws.Range("C11").Value = 0.03175288002049115
$ws.Range("D11").Value = 0.1946671942227312
$ws.Range("E11").Value = 0.18194385639503
$ws.Range("F11").Value = 1.626325773391898
$ws.Range("G11").Value = 0.9862726151528136
$ws.Range("H11").Value = 0.9799738081788405
$ws.Range("I11").Value = 0.8597837055919797
$ws.Range("J11").Value = 0.2140722647758508
$ws.Range("K11").Value = 3.05565635819238
$ws.Range("C12").Value = 0.03225217591995033
$ws.Range("D12").Value = 0.1965139665572337
$ws.Range("E12").Value = 0.1833936881828677
$ws.Range("F12").Value = 1.626754908650383
$ws.Range("G12").Value = 0.9860831862395543
$ws.Range("H12").Value = 0.9780551964002058
$ws.Range("I12").Value = 0.8598651689140695
$ws.Range("J12").Value = 0.2155654057072383
$ws.Range("K12").Value = 3.121282410437914
$ws.Range("C13").Value = 0.03214468284807026
$ws.Range("D13").Value = 0.196116025634538
$ws.Range("E13").Value = 0.1830809456276867
$ws.Range("F13").Value = 1.626653153969812
$ws.Range("G13").Value = 0.9861166056288937
$ws.Range("H13").Value = 0.9784633389128743
$ws.Range("I13").Value = 0.8598420217217537
$ws.Range("J13").Value = 0.2152429735465944
$ws.Range("K13").Value = 3.107148704346287
$ws.Range("C14").Value = 0.03179397504840153
$ws.Range("D14").Value = 0.1948190365951774
$ws.Range("E14").Value = 0.1820629117233565
$ws.Range("F14").Value = 1.626356874839857
$ws.Range("G14").Value = 0.9862537082253056
$ws.Range("H14").Value = 0.9798136812329545
$ws.Range("I14").Value = 0.8597878842214186
$ws.Range("J14").Value = 0.2141947200965859
$ws.Range("K14").Value = 3.06105546788325
$ws.Range("C15").Value = 0.03157904181193771
$ws.Range("D15").Value = 0.1940251963127224
$ws.Range("E15").Value = 0.1814407869251227
$ws.Range("F15").Value = 1.626202702264905
$ws.Range("G15").Value = 0.9863592668602337
$ws.Range("H15").Value = 0.9806556269148956
$ws.Range("I15").Value = 0.8597711140470281
$ws.Range("J15").Value = 0.21355514389316
$ws.Range("K15").Value = 3.032821945788612
$ws.Range("C16").Value = 0.03034560330341662
$ws.Range("D16").Value = 0.1894853198465256
$ws.Range("E16").Value = 0.1778978691538029
$ws.Range("F16").Value = 1.625736842323036
$ws.Range("G16").Value = 0.9872939222638877
$ws.Range("H16").Value = 0.9857072825976587
$ws.Range("I16").Value = 0.8599255943130046
$ws.Range("J16").Value = 0.2099283773250846
$ws.Range("K16").Value = 2.87102703774184
$ws.Range("C17").Value = 0.02958738796201033
$ws.Range("D17").Value = 0.1867086177642676
$ws.Range("E17").Value = 0.1757442362398862
$ws.Range("F17").Value = 1.625822743191321
$ws.Range("G17").Value = 0.9881604393554682
$ws.Range("H17").Value = 0.9890082843149202
$ws.Range("I17").Value = 0.860243267890425
$ws.Range("J17").Value = 0.2077376452967172
$ws.Range("K17").Value = 2.771772503853356
$ws.Range("C18").Value = 0.02915072217712122
$ws.Range("D18").Value = 0.1851146324726898
$ws.Range("E18").Value = 0.1745128112387206
$ws.Range("F18").Value = 1.626008359249667
$ws.Range("G18").Value = 0.9887661287449845
$ws.Range("H18").Value = 0.9909809967403334
$ws.Range("I18").Value = 0.8605076268608727
$ws.Range("J18").Value = 0.20649012031474
$ws.Range("K18").Value = 2.714685002091528
$ws.Range("C19").Value = 0.02900277896498693
$ws.Range("D19").Value = 0.1845754714268679
$ws.Range("E19").Value = 0.1740971235126736
$ws.Range("F19").Value = 1.626094562863983
$ws.Range("G19").Value = 0.9889895909841044
$ws.Range("H19").Value = 0.9916616284591271
$ws.Range("I19").Value = 0.8606111301817165
$ws.Range("J19").Value = 0.2060698781579617
$ws.Range("K19").Value = 2.695356394950068
$ws.Range("C20").Value = 0.02966815941034184
$ws.Range("D20").Value = 0.1870038821629976
$ws.Range("E20").Value = 0.1759727398787803
$ws.Range("F20").Value = 1.625799493188069
$ws.Range("G20").Value = 0.988057083076157
$ws.Range("H20").Value = 0.9886492181395852
$ws.Range("I20").Value = 0.8602009952253979
$ws.Range("J20").Value = 0.20796955533622
$ws.Range("K20").Value = 2.782338219415067
$ws.Range("C21").Value = 0.03189701034116865
$ws.Range("D21").Value = 0.1951998682025362
$ws.Range("E21").Value = 0.1823616306984874
$ws.Range("F21").Value = 1.626438206299696
$ws.Range("G21").Value = 0.9862089383320551
$ws.Range("H21").Value = 0.9794139631628127
$ws.Range("I21").Value = 0.8598003683952413
$ws.Range("J21").Value = 0.2145020944292213
$ws.Range("K21").Value = 3.074594189258107
$ws.Range("C22").Value = 0.03334858042740052
$ws.Range("D22").Value = 0.200583475635483
$ws.Range("E22").Value = 0.1866020530123151
$ws.Range("F22").Value = 1.628076794450507
$ws.Range("G22").Value = 0.9859657381291953
$ws.Range("H22").Value = 0.9740411539088853
$ws.Range("I22").Value = 0.8602713992908946
$ws.Range("J22").Value = 0.2188837107627961
$ws.Range("K22").Value = 3.265598117200284
$ws.Range("C23").Value = 0.03257432379446357
$ws.Range("D23").Value = 0.197707693616664
$ws.Range("E23").Value = 0.1843329192172405
$ws.Range("F23").Value = 1.627090113927878
$ws.Range("G23").Value = 0.9860068228182826
$ws.Range("H23").Value = 0.9768478951936999
$ws.Range("I23").Value = 0.8599526601381626
$ws.Range("J23").Value = 0.2165348578975426
$ws.Range("K23").Value = 3.163656493406904
$ws.Range("C24").Value = 0.0296316449968046
$ws.Range("D24").Value = 0.1868703857168583
$ws.Range("E24").Value = 0.1758694124304441
$ws.Range("F24").Value = 1.625809580260849
$ws.Range("G24").Value = 0.9881034756478329
$ws.Range("H24").Value = 0.9888113186162144
$ws.Range("I24").Value = 0.8602198521992221
$ws.Range("J24").Value = 0.2078646715750523
$ws.Range("K24").Value = 2.777561536370797
$ws.Range("C25").Value = 0.02643613307540704
$ws.Range("D25").Value = 0.175299425842951
$ws.Range("E25").Value = 0.1670187357216406
$ws.Range("F25").Value = 1.629630877107573
$ws.Range("G25").Value = 0.9944615000302406
$ws.Range("H25").Value = 1.004545382183451
$ws.Range("I25").Value = 0.8636258686906473
$ws.Range("J25").Value = 0.1989909355026356
$ws.Range("K25").Value = 2.36114304678182
